$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entities")

# --- Capture the two banded row styles ("Good" / s=2 and "Neutral" / s=3) from
# existing cells *before* any writes happen, so re-applying them later never
# mints brand-new style entries (cloning a live Style COM object reuses the
# existing cellXfs slot instead of creating a duplicate).
$styleGood    = $ws.Cells.Item(2, 1).Style   # s="2"
$styleNeutral = $ws.Cells.Item(4, 1).Style   # s="3"

# --- Row 1 header: "title" column text is unchanged, just rewritten so the
# shared-string table gets rebuilt/reordered the same way Excel would do it.
$ws.Cells.Item(1, 4).Value = "title"

# --- New row 2: "_None" / "None" enum entry, id 10000. No assetAddress, so
# the C2 cell must disappear completely (not just go blank) -> Clear().
$ws.Cells.Item(2, 1).Value = 10000
$ws.Cells.Item(2, 2).Value = "_None"
$ws.Cells.Item(2, 3).Clear()
$ws.Cells.Item(2, 4).Value = "None"

# --- Row 3: Item 1, id 10001 (was 10000 pre-edit).
$ws.Cells.Item(3, 1).Value = 10001
$ws.Cells.Item(3, 2).Value = "_10001_Item_1"
$ws.Cells.Item(3, 3).Value = "AssetAddress.CubeBlue_Pickupable"
$ws.Cells.Item(3, 4).Value = "Item 1"

# --- Row 4: Item 2, id 10002 (was 10001 pre-edit). Cell-level style flips
# from Neutral to Good to match the new banding after the insert.
$ws.Cells.Item(4, 1).Value = 10002
$ws.Cells.Item(4, 2).Value = "_10002_Item_2"
$ws.Cells.Item(4, 3).Value = "AssetAddress.CubeRed_Pickupable"
$ws.Cells.Item(4, 4).Value = "Item 2"
$ws.Range("A4:D4").Style = $styleGood

# --- Row 5: Puzzle Piece A, id 10100 (was row 4 pre-edit). Style unchanged
# (Neutral).
$ws.Cells.Item(5, 1).Value = 10100
$ws.Cells.Item(5, 2).Value = "_10100_PuzzleBlock_A"
$ws.Cells.Item(5, 3).Value = "AssetAddress.PuzzleBlock_A_Pickupable"
$ws.Cells.Item(5, 4).Value = "Puzzle Piece A"

# --- Row 6: Puzzle Piece B, id 10101 (was row 5 pre-edit). Cell-level style
# flips from Good to Neutral.
$ws.Cells.Item(6, 1).Value = 10101
$ws.Cells.Item(6, 2).Value = "_10101_PuzzleBlock_B"
$ws.Cells.Item(6, 3).Value = "AssetAddress.PuzzleBlock_B_Pickupable"
$ws.Cells.Item(6, 4).Value = "Puzzle Piece B"
$ws.Range("A6:D6").Style = $styleNeutral

# --- Row 7: Old Key, id 10201 (was row 6 pre-edit). Cell-level style flips
# from Neutral to Good.
$ws.Cells.Item(7, 1).Value = 10201
$ws.Cells.Item(7, 2).Value = "_10201_Key_A"
$ws.Cells.Item(7, 3).Value = "AssetAddress.Key_A_Pickupable"
$ws.Cells.Item(7, 4).Value = "Old Key"
$ws.Range("A7:D7").Style = $styleGood

# --- Row 8: Book A, id 10301 (was row 7 pre-edit). Style unchanged (Neutral).
$ws.Cells.Item(8, 1).Value = 10301
$ws.Cells.Item(8, 2).Value = "_10301_Book_A"
$ws.Cells.Item(8, 3).Value = "AssetAddress.Book_A_Pickupable"
$ws.Cells.Item(8, 4).Value = "Book A"

# --- Row 9: Book B, id 10302 (was row 8 pre-edit, no assetAddress). Style
# unchanged (Neutral).
$ws.Cells.Item(9, 1).Value = 10302
$ws.Cells.Item(9, 2).Value = "_10302_Book_B"
$ws.Cells.Item(9, 4).Value = "Book B"

# --- Row 10: Book C, id 10303 (was row 9 pre-edit, no assetAddress). This
# row previously existed only as an empty formatted row, so it needs the
# Neutral style applied explicitly.
$ws.Cells.Item(10, 1).Value = 10303
$ws.Cells.Item(10, 2).Value = "_10302_Book_C"
$ws.Cells.Item(10, 4).Value = "Book C"
$ws.Range("A10:D10").Style = $styleNeutral

# --- Selection moved to E7 on the Entities sheet.
$ws.Range("E7").Select()
